$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: header row for the new "New User Registration" test case ---
$ws.Range("A7").Value = "TC-003"
$ws.Range("B7").Value = "expTitle"
$ws.Range("C7").Value = "username"
$ws.Range("D7").Value = "password"
$ws.Range("E7").Value = "confirm Password"
$ws.Range("F7").Value = "Full  Name"
$ws.Range("G7").Value = "email id"
$ws.Range("H7").Value = "captcha"

# --- Row 8: data row for the new test case ---
$ws.Range("A8").Value = "TC-003"

# B8 reuses the same "expTitle answer" look (Comic Sans MS 14, custom blue/purple)
# that is already applied to D2 / D5, so copy that formatting across before
# setting the value.
[void]$ws.Range("D2").Copy()
[void]$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "Adactin.com - New User Registration"

$ws.Range("C8").Value = "kiran2403"
$ws.Range("D8").Value = "kiran123"
$ws.Range("E8").Value = "kiran123"
$ws.Range("F8").Value = "Kiran Kumar"
$ws.Range("G8").Value = "abcd@gmail.com"
$ws.Range("H8").Value = "Hyderabad"

# Hyperlink on the captcha/email cell G8, pointing at the registration page.
[void]$ws.Hyperlinks.Add($ws.Range("G8"), "http://www.adactin.com/HotelApp/index.php")

# --- Column widths: widen B for the longer title text, and size the four
# new columns E-H to fit their content (matches Excel's own best-fit sizing
# as closely as this engine's 1/6-character width grid allows). ---
$ws.Columns("B:B").ColumnWidth = 27.608072916666668
$ws.Columns("E:E").ColumnWidth = 25.166666666666668
$ws.Columns("F:F").ColumnWidth = 14.944010416666666
$ws.Columns("G:G").ColumnWidth = 14.385416666666666
$ws.Columns("H:H").ColumnWidth = 15.385416666666666

# Move the active selection like the saved workbook shows.
[void]$ws.Range("B12").Select()
